# Fruta / hortaliza, semanal
# Insert 4 new weekly records (week of 2021-11-27 / serial 44509) at the
# top of the "Frutilla" price table, pushing the existing historical rows
# (658..689) down to (662..693).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the first data row of the table (row 658),
# shifting everything below (rows 658:689) down to rows 662:693.
$ws.Rows("658:661").Insert()

# Columns that are constant for every "Frutilla" record in this sheet.
$const = @{
    A = 6
    B = "Mercado Mayorista Lo Valledor de Santiago"
    C = "Metropolitana"
    E = 13
    F = "Fruta"
    G = 100101
    H = "Berries"
    I = 100112025
    J = "Frutilla"
    K = "Sin especificar"
    Q = '$/bandeja 7 kilos'
    T = 7
}

# New rows: Fecha(D), Calidad(L), Volumen(M), Precio minimo(N),
# Precio maximo(O), Precio promedio ponderado(P), Origen(R), Precio $/Kg(S)
$newRows = @(
    @{ Row = 658; D = 44509; L = "Especial"; M = 2200; N = 6000; O = 7000; P = 6420; R = "Provincia de Melipilla"; S = 917 },
    @{ Row = 659; D = 44509; L = "Primera";  M = 1200; N = 5000; O = 5500; P = 5188; R = "Provincia de Melipilla"; S = 741 },
    @{ Row = 660; D = 44509; L = "Segunda";  M = 1450; N = 3500; O = 4500; P = 3879; R = "Provincia de Melipilla"; S = 554 },
    @{ Row = 661; D = 44509; L = "Tercera";  M = 400;  N = 2500; O = 3000; P = 2750; R = "Provincia de Melipilla"; S = 393 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    foreach ($col in $const.Keys) {
        $ws.Cells.Item($row, [int][char]$col - [int][char]'A' + 1).Value = $const[$col]
    }

    $ws.Cells.Item($row, 4).Value = $r.D   # D
    $ws.Cells.Item($row, 12).Value = $r.L  # L
    $ws.Cells.Item($row, 13).Value = $r.M  # M
    $ws.Cells.Item($row, 14).Value = $r.N  # N
    $ws.Cells.Item($row, 15).Value = $r.O  # O
    $ws.Cells.Item($row, 16).Value = $r.P  # P
    $ws.Cells.Item($row, 18).Value = $r.R  # R
    $ws.Cells.Item($row, 19).Value = $r.S  # S
}
